# Suite.xlsx edit: swap the Runmode flags for the two test suites and
# move the active selection from B3 to B2.
#
# Before:  B2 (BankManagerSuite) = "Y", B3 (CustomerSuite) = "N"
# After:   B2 (BankManagerSuite) = "N", B3 (CustomerSuite) = "Y"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "N"
$ws.Range("B3").Value = "Y"

# Selection moves to B2 (was B3).
$ws.Range("B2").Select()
